$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '38.820.96'
$ws.Range("E2").Value = '  +2.77%  '
$ws.Range("D3").Value = '2.091.90'
$ws.Range("E3").Value = '  +2.34%  '
$ws.Range("E4").Value = '  -0.11%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '228.80'
$ws.Range("E5").Value = '  +0.56%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.615'
$ws.Range("E6").Value = '  +0.91%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '60.42'
$ws.Range("E7").Value = '  +0.57%  '
$ws.Range("E8").Value = '  -0.09%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.385'
$ws.Range("E9").Value = '  +2.39%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0836'
$ws.Range("E10").Value = '  +0.13%  '
$ws.Range("E11").Value = '  -0.38%  '
$ws.Range("D12").Value = '2.402.43'
$ws.Range("E12").Value = '  +2.29%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '14.98'
$ws.Range("E13").Value = '  +4.50%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '21.84'
$ws.Range("E14").Value = '  +2.14%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.796'
$ws.Range("E15").Value = '  +4.41%  '
$ws.Range("E16").Value = '  -0.09%  '
$ws.Range("D17").Value = '2.094.74'
$ws.Range("E17").Value = '  +2.43%  '
$ws.Range("D18").Value = '38.689.27'
$ws.Range("E18").Value = '  +2.56%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '71.63'
$ws.Range("E19").Value = '  +3.32%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.04'
$ws.Range("E20").Value = '  +2.21%  '
$ws.Range("D21").Value = '0.0₃0837'
$ws.Range("E21").Value = '  +1.12%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '227.29'
$ws.Range("E22").Value = '  +2.13%  '
$ws.Range("E23").Value = '  -0.42%  '
$ws.Range("E24").Value = '  -0.02%  '
$ws.Range("E25").Value = '  +3.00%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '171.03'
$ws.Range("E26").Value = '  +1.08%  '
$ws.Range("E27").Value = '  +2.01%  '
$ws.Range("E28").Value = '  +9.13%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.46'
$ws.Range("E29").Value = '  +13.73%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '19.17'
$ws.Range("E30").Value = '  +2.22%  '
$ws.Range("E31").Value = '  +0.97%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.37'
$ws.Range("E32").Value = '  +5.34%  '
$ws.Range("E33").Value = '  +3.04%  '
$ws.Range("E34").Value = '  +4.67%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0612'
$ws.Range("E35").Value = '  +1.65%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.44'
$ws.Range("E36").Value = '  -1.51%  '
$ws.Range("E37").Value = '  +2.05%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.58'
$ws.Range("E38").Value = '  +3.67%  '
$ws.Range("E39").Value = '  -0.11%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '18.16'
$ws.Range("E40").Value = '  +0.80%  '
$ws.Range("D41").Value = '1.541.12'
$ws.Range("E41").Value = '  +0.65%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '100.88'
$ws.Range("E42").Value = '  +3.38%  '
$ws.Range("E43").Value = '  +4.29%  '
$ws.Range("E44").Value = '  -0.82%  '
$ws.Range("E45").Value = '  +3.59%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '7.63'
$ws.Range("E46").Value = '  +8.31%  '
$ws.Range("E47").Value = '  +1.06%  '
$ws.Range("E48").Value = '  -1.25%  '
$ws.Range("E49").Value = '  +2.80%  '
$ws.Range("E50").Value = '  +1.08%  '
$ws.Range("D51").Value = '2.289.07'
$ws.Range("E51").Value = '  +2.35%  '
